$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("C2").Value = 28

# Row 3
$ws.Range("B3").Value = "<dumber>"
$ws.Range("C3").Value = 28

# Row 4
$ws.Range("C4").Value = 27

# Row 5
$ws.Range("B5").Value = "<fot>"
$ws.Range("C5").Value = 30

# Row 6
$ws.Range("C6").Value = 19

# Row 8
$ws.Range("B8").Value = "<ostar>"
$ws.Range("C8").Value = 24

# Row 9
$ws.Range("C9").Value = 35

# Row 10
$ws.Range("B10").Value = "<eis>"
$ws.Range("C10").Value = 24

# Row 11
$ws.Range("C11").Value = 26

# Row 12
$ws.Range("C12").Value = 33

# Row 13
$ws.Range("B13").Value = "<vite>"
$ws.Range("C13").Value = 32

# Row 14
$ws.Range("B14").Value = "<alt>"

# Row 15
$ws.Range("C15").Value = 23

# Row 16
$ws.Range("C16").Value = 23

# Row 17
$ws.Range("C17").Value = 27

# Row 18
$ws.Range("B18").Value = "<his>"
$ws.Range("C18").Value = 21
